# The workbook previously had an unused/empty column H sitting between the
# "e_yy_std" block (B:G) and the "Timestamp"/"elapsed time" block that lived
# in columns I:J. This edit removes that empty spacer column, which shifts
# the Timestamp/elapsed-time data left from I:J to H:I (renamed data sets /
# graphical time-offset columns per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty column H - shifts old I -> H (Timestamp) and old J -> I
# (elapsed time), matching the new H:I layout.
$ws.Columns("H").Delete()

# Re-fit the surviving columns to their (now slightly different) content
# widths, closest to how Excel recalculated them after the shift.
$ws.Columns("C").ColumnWidth = 27.5
$ws.Columns("H").ColumnWidth = 18.33
$ws.Columns("I").ColumnWidth = 11.17

# Restore the user's on-screen selection to a neutral cell after the edit.
$ws.Range("L3").Select()
